$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.007.47'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '3.395.48'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.89'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.75%  '
$ws.Range('D7').Value = '3.397.48'
$ws.Range('E7').Value = '  -1.27%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.52'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.395'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').Value = '3.970.19'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.50'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '3.391.87'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '61.044.19'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.23'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.04'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '385.93'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.560'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '74.10'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000118'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.08%  '
$ws.Range('D27').Value = '3.528.10'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.42'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.03'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('E33').Value = '  -2.37%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.71'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.02'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '165.70'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').Value = '3.422.35'
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.01'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.96%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.49'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.87%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '28.28'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0777'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.781'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.57%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.19'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.43'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E47').Value = '  -3.10%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.14'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('D49').Value = '2.488.16'
$ws.Range('E49').Value = '  -4.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.48'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.83'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.95%  '
